# Update the "想去人数" (want-to-go count) figures in the F column
# across the relevant worksheets, per the latest data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1644
$ws1.Range("F3").Value = 9043
$ws1.Range("F4").Value = 109
$ws1.Range("F5").Value = 502
$ws1.Range("F6").Value = 693
$ws1.Range("F7").Value = 693
$ws1.Range("F8").Value = 189
$ws1.Range("F9").Value = 50
$ws1.Range("F10").Value = 83
$ws1.Range("F11").Value = 5372
$ws1.Range("F12").Value = 59
$ws1.Range("F15").Value = 4336
$ws1.Range("F18").Value = 1146
$ws1.Range("F19").Value = 15
$ws1.Range("F20").Value = 335
$ws1.Range("F21").Value = 14
$ws1.Range("F22").Value = 251
$ws1.Range("F24").Value = 2677
$ws1.Range("F25").Value = 120

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 4
$ws2.Range("F3").Value = 39

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1644
$ws4.Range("F3").Value = 9043
$ws4.Range("F4").Value = 109
$ws4.Range("F5").Value = 4
$ws4.Range("F6").Value = 502
$ws4.Range("F7").Value = 693
$ws4.Range("F8").Value = 693
$ws4.Range("F9").Value = 189
$ws4.Range("F10").Value = 50
$ws4.Range("F11").Value = 83
$ws4.Range("F12").Value = 5374
$ws4.Range("F13").Value = 59
$ws4.Range("F16").Value = 4336
$ws4.Range("F19").Value = 1146
$ws4.Range("F20").Value = 15
$ws4.Range("F21").Value = 335
$ws4.Range("F22").Value = 14
$ws4.Range("F23").Value = 251
$ws4.Range("F25").Value = 2677
$ws4.Range("F26").Value = 39
$ws4.Range("F27").Value = 120
